# #3 added test harness for login and added login fixed.
#
# Column A held bare usernames paired with the *wrong* password in column B
# (rows were a shuffled lookup table). This rewrites column A as the real
# sign-in email (user@TestIncidentQueue.onmicrosoft.com) turned into a
# mailto: hyperlink, realigns column B so every row's password actually
# belongs to that row's account, and adds a fixed P@ssw0rd1 test login
# (hyperlinked to the sign-in page) in row 1 next to RickG.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: realign every row's password to its own account first (these
# reuse the six passwords already in the shared-string table, just in their
# new row positions, so the table's surviving entries keep their order).
$ws.Range("B2").Value = "Labo0749"
$ws.Range("B3").Value = "Zuwu5875"
$ws.Range("B4").Value = "Pufa7292"
$ws.Range("B5").Value = "Hoyo4800"
$ws.Range("B6").Value = "Puva8501"
$ws.Range("B7").Value = "Tuxo4459"

# --- Column A: sign-in emails, written in this order so new shared-string
# entries are appended as RickG, DarylD, CarlG, GlennR, MaggieG, HershelG,
# Michonne - then B1's fixed test password last.
$ws.Range("A1").Value = "RickG@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A2").Value = "DarylD@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A4").Value = "CarlG@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A5").Value = "GlennR@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A6").Value = "MaggieG@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A7").Value = "HershelG@TestIncidentQueue.onmicrosoft.com"
$ws.Range("A3").Value = "Michonne@TestIncidentQueue.onmicrosoft.com"
$ws.Range("B1").Value = "P@ssw0rd1"

# Hyperlink the account emails (added in this order so the relationship ids
# land A1,A2,A4,A5,A6,A7,A3 then the B1 login link last).
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:RickG@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:DarylD@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:CarlG@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:GlennR@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:MaggieG@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:HershelG@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:Michonne@TestIncidentQueue.onmicrosoft.com") | Out-Null

# The fixed test login harness: password cell links to the sign-in page.
$ws.Hyperlinks.Add($ws.Range("B1"), "https://login.microsoftonline.com/") | Out-Null

# Column A needs to be wide enough to show the full email address.
$ws.Columns("A").ColumnWidth = 41.75

# Leave the selection where the editor ended up.
$ws.Range("D4").Select() | Out-Null
